$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P3").Value = 1.99
$ws.Range("Q3").Value = 1.75
$ws.Range("H4").Value = 2.16
$ws.Range("I4").Value = 2.44
$ws.Range("J4").Value = 3.4
$ws.Range("Q4").Value = 1.84
$ws.Range("R4").Value = 1.32
$ws.Range("V4").Value = 1.63
$ws.Range("AG4").Value = 20
$ws.Range("AI4").Value = 48
$ws.Range("AJ4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("I5").Value = 4.5
$ws.Range("L5").Value = 1.31
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 3.6
$ws.Range("O5").Value = 1.3
$ws.Range("R5").Value = 1.34
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 1.73
$ws.Range("U5").Value = 2.08
$ws.Range("V5").Value = 1.29
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 17.5
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 19.5
$ws.Range("AE5").Value = 980
$ws.Range("AF5").Value = 16.5
$ws.Range("AG5").Value = 13
$ws.Range("AI5").Value = 65
$ws.Range("AL5").Value = 980
$ws.Range("AN5").Value = 980
$ws.Range("F6").Value = 8.6
$ws.Range("G6").Value = 10.5
$ws.Range("H6").Value = 1.33
$ws.Range("I6").Value = 1.4
$ws.Range("J6").Value = 5.7
$ws.Range("K6").Value = 6.8
$ws.Range("N6").Value = 6.4
$ws.Range("P6").Value = 2.98
$ws.Range("Q6").Value = 1.4
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 2.02
$ws.Range("V6").Value = 3.5
$ws.Range("W6").Value = 1.1
$ws.Range("Y6").Value = 1000
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 13
$ws.Range("AC6").Value = 18
$ws.Range("AD6").Value = 12
$ws.Range("AF6").Value = 110
$ws.Range("AJ6").Value = 310
$ws.Range("AK6").Value = 130
$ws.Range("AL6").Value = 100
$ws.Range("AM6").Value = 110
$ws.Range("AO6").Value = 4.2
$ws.Range("F7").Value = 2.16
$ws.Range("K7").Value = 3.6
$ws.Range("N7").Value = 1.1
$ws.Range("Q7").Value = 2.02
$ws.Range("S7").Value = 1.05
$ws.Range("AI7").Value = 85
$ws.Range("G8").Value = 1.34
$ws.Range("I8").Value = 1000
$ws.Range("V8").Value = 1.05
$ws.Range("W8").Value = 3.9
$ws.Range("U9").Value = 1.94
$ws.Range("H10").Value = 1.67
$ws.Range("I10").Value = 1.68
$ws.Range("Q10").Value = 1.85
$ws.Range("V10").Value = 2.46
$ws.Range("AF10").Value = 46
$ws.Range("AG10").Value = 22
$ws.Range("AH10").Value = 21
$ws.Range("AI10").Value = 34
$ws.Range("AG11").Value = 9.800000000000001
$ws.Range("I12").Value = 2.28
$ws.Range("M12").Value = 1.06
$ws.Range("T12").Value = 1.68
$ws.Range("G13").Value = 1.5
$ws.Range("W13").Value = 3
$ws.Range("F14").Value = 2.48
$ws.Range("G14").Value = 2.8
$ws.Range("S14").Value = 3.7
$ws.Range("T14").Value = 1.79
$ws.Range("AL14").Value = 55
$ws.Range("Q15").Value = 1.61
$ws.Range("F18").Value = 1.91
$ws.Range("N18").Value = 3.95
$ws.Range("R18").Value = 1.38
$ws.Range("S18").Value = 3.05
$ws.Range("T18").Value = 1.75
$ws.Range("U18").Value = 2.12
$ws.Range("P19").Value = 2.24
$ws.Range("Q19").Value = 1.15
$ws.Range("R19").Value = 1.61
$ws.Range("S19").Value = 2.02
$ws.Range("X19").Value = 34
$ws.Range("Z19").Value = 1000
$ws.Range("AI19").Value = 1000
$ws.Range("F20").Value = 1.44
$ws.Range("G20").Value = 1.63
$ws.Range("H20").Value = 6.4
$ws.Range("K20").Value = 5.6
$ws.Range("U20").Value = 1.7
$ws.Range("H22").Value = 2.38
$ws.Range("J22").Value = 3.4
$ws.Range("K22").Value = 4.1
$ws.Range("N23").Value = 4.4
$ws.Range("T23").Value = 1.76
$ws.Range("U23").Value = 2.04
$ws.Range("X23").Value = 21
$ws.Range("Y23").Value = 25
$ws.Range("Z23").Value = 55
$ws.Range("AB23").Value = 10.5
$ws.Range("AC23").Value = 11
$ws.Range("AD23").Value = 25
$ws.Range("AE23").Value = 85
$ws.Range("AF23").Value = 11
$ws.Range("AG23").Value = 10.5
$ws.Range("AH23").Value = 980
$ws.Range("AI23").Value = 75
$ws.Range("AJ23").Value = 16.5
$ws.Range("AK23").Value = 17
$ws.Range("AL23").Value = 980
$ws.Range("AN23").Value = 8
$ws.Range("AO23").Value = 110
$ws.Range("F24").Value = 5.3
$ws.Range("G24").Value = 7.4
$ws.Range("H24").Value = 1.53
$ws.Range("I24").Value = 1.69
$ws.Range("J24").Value = 4.6
$ws.Range("N24").Value = 2.56
$ws.Range("P24").Value = 2.56
$ws.Range("V24").Value = 2.46
$ws.Range("I25").Value = 3.65
$ws.Range("J25").Value = 2.92
$ws.Range("Q25").Value = 1.98
$ws.Range("R25").Value = 1.26
$ws.Range("V25").Value = 1.38
$ws.Range("F26").Value = 1.48
$ws.Range("G26").Value = 1.6
$ws.Range("H26").Value = 6.2
$ws.Range("I26").Value = 7.6
$ws.Range("J26").Value = 4.1
$ws.Range("K26").Value = 5.5
$ws.Range("T26").Value = 1.54
$ws.Range("U26").Value = 1.87
$ws.Range("W26").Value = 2.66
$ws.Range("X26").Value = 32
$ws.Range("Y26").Value = 40
$ws.Range("AB26").Value = 16
$ws.Range("AC26").Value = 17
$ws.Range("AD26").Value = 38
$ws.Range("AF26").Value = 15.5
$ws.Range("AG26").Value = 15
$ws.Range("AH26").Value = 30
$ws.Range("AJ26").Value = 21
$ws.Range("AK26").Value = 22
$ws.Range("AL26").Value = 44
$ws.Range("AN26").Value = 8.800000000000001
$ws.Range("F27").Value = 2.14
$ws.Range("H27").Value = 2.9
$ws.Range("I27").Value = 3.7
$ws.Range("J27").Value = 3.1
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 1.21
$ws.Range("N27").Value = 2.16
$ws.Range("P27").Value = 2.16
$ws.Range("Q27").Value = 1.48
$ws.Range("R27").Value = 1.48
$ws.Range("S27").Value = 2.16
$ws.Range("T27").Value = 1.04
$ws.Range("U27").Value = 1.04
$ws.Range("X27").Value = 990
$ws.Range("Y27").Value = 1000
$ws.Range("Z27").Value = 1000
$ws.Range("AA27").Value = 1000
$ws.Range("AB27").Value = 1000
$ws.Range("AC27").Value = 1000
$ws.Range("AD27").Value = 1000
$ws.Range("AE27").Value = 1000
$ws.Range("AF27").Value = 1000
$ws.Range("AG27").Value = 1000
$ws.Range("AH27").Value = 1000
$ws.Range("AI27").Value = 1000
$ws.Range("AJ27").Value = 1000
$ws.Range("AK27").Value = 1000
$ws.Range("AL27").Value = 1000
$ws.Range("AM27").Value = 1000
$ws.Range("AN27").Value = 1000
$ws.Range("AO27").Value = 1000
$ws.Range("N28").Value = 1.1
$ws.Range("P28").Value = 1.63
$ws.Range("Q28").Value = 1.99
$ws.Range("S28").Value = 2.52
$ws.Range("Y28").Value = 22
$ws.Range("Z28").Value = 50
$ws.Range("AC28").Value = 12
$ws.Range("AD28").Value = 28
$ws.Range("AE28").Value = 1000
$ws.Range("AF28").Value = 17
$ws.Range("AG28").Value = 15
$ws.Range("AI28").Value = 1000
$ws.Range("AJ28").Value = 32
$ws.Range("AK28").Value = 30
$ws.Range("AL28").Value = 1000
$ws.Range("AN28").Value = 22
$ws.Range("G30").Value = 4.4
$ws.Range("J30").Value = 3.8
$ws.Range("W30").Value = 1.3
$ws.Range("U31").Value = 1.78
$ws.Range("AK31").Value = 16.5
$ws.Range("AM31").Value = 170
$ws.Range("AA32").Value = 330
$ws.Range("K33").Value = 4.9
